$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-23 12:46:36"

# Update the "取得日時" (fetched-at) timestamp for existing rows 2-11 to the
# new run time.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Append the new case picked up by this run as row 12.
$row = 12
$ws.Cells.Item($row, 1).Value = $newTimestamp
$ws.Cells.Item($row, 2).Value = "【急募】JotformとGoogleスプレッドシート連携のエキスパート募集!"
$ws.Cells.Item($row, 3).Value = "システム開発"
$ws.Cells.Item($row, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item($row, 5).Value = "期限情報なし"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 6), "https://www.lancers.jp/work/detail/5395809") | Out-Null
$ws.Cells.Item($row, 6).Style = "Hyperlink"
$ws.Cells.Item($row, 7).Value = 10
